$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert two new rows at 18-19 (shifting the existing rows 18-77 down to
# 20-79), then populate them with the two new case records (7671 and 7673).
# ---------------------------------------------------------------------------
$ws.Range("A18:A19").EntireRow.Insert()

# --- New row 18 : Caso 7671 -------------------------------------------------
$ws.Cells.Item(18, 1).Value  = "'7671"
$ws.Cells.Item(18, 2).Value  = "'5/5/2025"
$ws.Cells.Item(18, 3).Value  = "MATIENZO, BENJAMIN, TENIENTE 1520"
$ws.Cells.Item(18, 4).Value  = "'14"
$ws.Cells.Item(18, 5).Value  = "'805655369"
$ws.Cells.Item(18, 6).Value  = "NEW"
$ws.Cells.Item(18, 7).Value  = "Pendiente"
$ws.Cells.Item(18, 8).Value  = "Picada"
$ws.Cells.Item(18, 9).Value  = 1
$ws.Cells.Item(18, 10).Value = "Cambio"
$ws.Cells.Item(18, 11).Value = "Sin equipos"
$ws.Cells.Item(18, 12).Value = "Terminal"
$ws.Cells.Item(18, 13).Value = -58.432419
$ws.Cells.Item(18, 14).Value = -34.566431
$ws.Cells.Item(18, 15).Value = "Palermo"
$ws.Cells.Item(18, 16).Value = "Capital Sur"
$ws.Cells.Item(18, 17).Value = "BLO-I"
$ws.Cells.Item(18, 18).Value = "Fuera de Poligono OVL"

# --- New row 19 : Caso 7673 -------------------------------------------------
$ws.Cells.Item(19, 1).Value  = "'7673"
$ws.Cells.Item(19, 2).Value  = "'5/7/2025"
$ws.Cells.Item(19, 3).Value  = "CAMPOS, LUIS M. AV. 1336"
$ws.Cells.Item(19, 4).Value  = "'14"
$ws.Cells.Item(19, 5).Value  = "'805722772"
$ws.Cells.Item(19, 6).Value  = "NEW"
$ws.Cells.Item(19, 7).Value  = "Pendiente"
$ws.Cells.Item(19, 8).Value  = "Picada"
$ws.Cells.Item(19, 9).Value  = 1
$ws.Cells.Item(19, 10).Value = "Cambio"
$ws.Cells.Item(19, 11).Value = "Sin equipos"
$ws.Cells.Item(19, 12).Value = "Pasante"
$ws.Cells.Item(19, 13).Value = -58.44191
$ws.Cells.Item(19, 14).Value = -34.564245
$ws.Cells.Item(19, 15).Value = "Colegiales"
$ws.Cells.Item(19, 16).Value = "Capital Norte"
$ws.Cells.Item(19, 17).Value = "BLO-L"
$ws.Cells.Item(19, 18).Value = "Fuera de Poligono OVL"

# ---------------------------------------------------------------------------
# Append a new row 80 with the last new case record (7676). At this point the
# previously-last row (old 77) now lives at row 79, so row 80 is a fresh row.
# ---------------------------------------------------------------------------
$ws.Cells.Item(80, 1).Value  = "'7676"
$ws.Cells.Item(80, 2).Value  = "'10/28/2025"
$ws.Cells.Item(80, 3).Value  = "RIVAS, GRAL. 2365"
$ws.Cells.Item(80, 4).Value  = "'11"
$ws.Cells.Item(80, 5).Value  = "'810461116"
$ws.Cells.Item(80, 6).Value  = "NEW"
$ws.Cells.Item(80, 7).Value  = "Pendiente"
$ws.Cells.Item(80, 8).Value  = "Poste apoyado en arbol cambiar o desmontar"
$ws.Cells.Item(80, 9).Value  = 1
$ws.Cells.Item(80, 10).Value = "Cambio"
$ws.Cells.Item(80, 11).Value = "Sin equipos"
$ws.Cells.Item(80, 12).Value = "Poste"
$ws.Cells.Item(80, 13).Value = -58.482578
$ws.Cells.Item(80, 14).Value = -34.59884
$ws.Cells.Item(80, 15).Value = "Paternal"
$ws.Cells.Item(80, 16).Value = "Capital Norte"
$ws.Cells.Item(80, 17).Value = "PUE-A"
$ws.Cells.Item(80, 18).Value = "Fuera de Poligono OVL"
